$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.02
    "C2" = 1.065169162970301
    "D2" = 1.067213479391637
    "E2" = 1.077738578085972
    "F2" = 1.082623260000702
    "I2" = 1.043626918840186
    "J2" = 1.07012565635057
    "K2" = 1.069921945177952
    "L2" = 1.080419088038693
    "M2" = 1.085290995118362
    "N2" = 1.026879362773818
    "B3" = 1.02
    "C3" = 1.066647031909996
    "D3" = 1.068386248650872
    "E3" = 1.079165028368179
    "F3" = 1.084031719245346
    "I3" = 1.043943249152399
    "J3" = 1.07125652186067
    "K3" = 1.070909567378499
    "L3" = 1.081661804069418
    "M3" = 1.08651669651655
    "N3" = 1.027271049835224
    "B4" = 1.02
    "C4" = 1.067602337643052
    "D4" = 1.069143997951613
    "E4" = 1.080087383030484
    "F4" = 1.084942377958715
    "I4" = 1.044146084143764
    "J4" = 1.071986849688108
    "K4" = 1.071546936660056
    "L4" = 1.08246474038502
    "M4" = 1.087308563333756
    "N4" = 1.027523662580099
    "B5" = 1.02
    "C5" = 1.068003721110786
    "D5" = 1.06946229388834
    "E5" = 1.080474989592401
    "F5" = 1.085325054324697
    "I5" = 1.04423091405485
    "J5" = 1.072293544577684
    "K5" = 1.071814486339789
    "L5" = 1.082802016119561
    "M5" = 1.0876411709813
    "N5" = 1.027629662491107
    "B6" = 1.02
    "C6" = 1.068071102025419
    "D6" = 1.069515721938299
    "E6" = 1.080540061717111
    "F6" = 1.085389297867973
    "I6" = 1.044245131497016
    "J6" = 1.072345020488832
    "K6" = 1.071859385741306
    "L6" = 1.082858629995085
    "M6" = 1.087697000212143
    "N6" = 1.02764744874397
    "B7" = 1.02
    "C7" = 1.067607701833118
    "D7" = 1.069148252061432
    "E7" = 1.080092562836484
    "F7" = 1.084947491939338
    "I7" = 1.044147219379872
    "J7" = 1.071990949070716
    "K7" = 1.071550513240412
    "L7" = 1.08246924816905
    "M7" = 1.0873130087996
    "N7" = 1.027525079734273
    "B8" = 1.02
    "C8" = 1.065668821162009
    "D8" = 1.067610054712247
    "E8" = 1.078220792226955
    "F8" = 1.083099405183253
    "I8" = 1.043734208451894
    "J8" = 1.070508133912896
    "K8" = 1.070256068609202
    "L8" = 1.08083931778508
    "M8" = 1.085705487012321
    "N8" = 1.027011908987229
    "B9" = 1.02
    "C9" = 1.062244526114509
    "D9" = 1.06489086387414
    "E9" = 1.074917248388542
    "F9" = 1.079837184750526
    "I9" = 1.042992183139508
    "J9" = 1.06788416389505
    "K9" = 1.067961989329671
    "L9" = 1.077957889274508
    "M9" = 1.082863092454155
    "N9" = 1.026101178212415
    "B10" = 1.02
    "C10" = 1.059956069141445
    "D10" = 1.063071968772416
    "E10" = 1.072711027930265
    "F10" = 1.077658247124888
    "I10" = 1.042487829605995
    "J10" = 1.066127130169553
    "K10" = 1.066423551330736
    "L10" = 1.076030393676277
    "M10" = 1.080961323122769
    "N10" = 1.025489592374622
    "B11" = 1.02
    "C11" = 1.058963727990521
    "D11" = 1.062282861823243
    "E11" = 1.071754720659648
    "F11" = 1.076713692687348
    "I11" = 1.042267124745111
    "J11" = 1.065364420186706
    "K11" = 1.065755189130886
    "L11" = 1.075194147017711
    "M11" = 1.080136147706824
    "N11" = 1.025223697448806
    "B12" = 1.02
    "C12" = 1.058594906740055
    "D12" = 1.061989519897094
    "E12" = 1.071399348939148
    "F12" = 1.076362677540509
    "I12" = 1.042184795273188
    "J12" = 1.065080824285194
    "K12" = 1.065506593070953
    "L12" = 1.07488327720637
    "M12" = 1.079829380122
    "N12" = 1.02512476905163
    "B13" = 1.02
    "C13" = 1.05867403029894
    "D13" = 1.062052453338138
    "E13" = 1.071475584571679
    "F13" = 1.076437979075317
    "I13" = 1.04220247108595
    "J13" = 1.065141669847705
    "K13" = 1.065559933077975
    "L13" = 1.074949971270532
    "M13" = 1.079895194700446
    "N13" = 1.025145996941947
    "B14" = 1.02
    "C14" = 1.05893324567483
    "D14" = 1.062258618860097
    "E14" = 1.071725348755105
    "F14" = 1.076684681080086
    "I14" = 1.042260326512997
    "J14" = 1.06534098403902
    "K14" = 1.065734646997004
    "L14" = 1.075168455568082
    "M14" = 1.080110795554212
    "N14" = 1.025215523335935
    "B15" = 1.02
    "C15" = 1.059092927235804
    "D15" = 1.062385613284619
    "E15" = 1.071879215763014
    "F15" = 1.076836660197991
    "I15" = 1.042295926736019
    "J15" = 1.065463749302343
    "K15" = 1.065842249264692
    "L15" = 1.075303037573926
    "M15" = 1.0802435996379
    "N15" = 1.025258339164755
    "B16" = 1.02
    "C16" = 1.060021896902317
    "D16" = 1.063124306882058
    "E16" = 1.072774473287852
    "F16" = 1.07772091122013
    "I16" = 1.042502428099097
    "J16" = 1.066177708154072
    "K16" = 1.066467861392006
    "L16" = 1.076085857800451
    "M16" = 1.081016051056741
    "N16" = 1.025507216159794
    "B17" = 1.02
    "C17" = 1.060604228389394
    "D17" = 1.063587261096189
    "E17" = 1.073335772580632
    "F17" = 1.078275290273239
    "I17" = 1.042631339474225
    "J17" = 1.066625042319981
    "K17" = 1.066859696727048
    "L17" = 1.076576460119445
    "M17" = 1.081500131089885
    "N17" = 1.025663041262156
    "B18" = 1.02
    "C18" = 1.06094375559702
    "D18" = 1.063857148786213
    "E18" = 1.073663072989896
    "F18" = 1.078598548343681
    "I18" = 1.042706307919535
    "J18" = 1.06688578155918
    "K18" = 1.067088034896889
    "L18" = 1.076862463443307
    "M18" = 1.081782323464912
    "N18" = 1.025753827883322
    "B19" = 1.02
    "C19" = 1.0610595026094
    "D19" = 1.063949148959247
    "E19" = 1.073774657888241
    "F19" = 1.078708753905997
    "I19" = 1.042731832388378
    "J19" = 1.066974655993809
    "K19" = 1.067165856369845
    "L19" = 1.076959956783697
    "M19" = 1.081878516244538
    "N19" = 1.025784766269314
    "B20" = 1.02
    "C20" = 1.060541763946943
    "D20" = 1.063537605594839
    "E20" = 1.073275560452577
    "F20" = 1.07821582122828
    "I20" = 1.042617531622165
    "J20" = 1.066577066588469
    "K20" = 1.066817678544626
    "L20" = 1.076523839385776
    "M20" = 1.08144821081836
    "N20" = 1.025646333424462
    "B21" = 1.02
    "C21" = 1.058856919384868
    "D21" = 1.06219791471091
    "E21" = 1.071651803804158
    "F21" = 1.076612038110059
    "I21" = 1.042243299204417
    "J21" = 1.065282299068033
    "K21" = 1.065683207432404
    "L21" = 1.075104124352753
    "M21" = 1.080047313719523
    "N21" = 1.025195054064917
    "B22" = 1.02
    "C22" = 1.057796305045098
    "D22" = 1.061354250323968
    "E22" = 1.070629974146045
    "F22" = 1.075602714902509
    "I22" = 1.042005979328936
    "J22" = 1.064466538282341
    "K22" = 1.064967970669562
    "L22" = 1.07421004048619
    "M22" = 1.079165002924417
    "N22" = 1.024910371796917
    "B23" = 1.02
    "C23" = 1.058358681190452
    "D23" = 1.061801622223779
    "E23" = 1.071171753761238
    "F23" = 1.076137869450415
    "I23" = 1.042131979658578
    "J23" = 1.064899150485322
    "K23" = 1.065347317632425
    "L23" = 1.0746841509462
    "M23" = 1.07963287775225
    "N23" = 1.025061377432855
    "B24" = 1.02
    "C24" = 1.06056998936305
    "D24" = 1.063560043228949
    "E24" = 1.073302768018343
    "F24" = 1.078242693044957
    "I24" = 1.042623771486508
    "J24" = 1.066598745322846
    "K24" = 1.066836665410277
    "L24" = 1.076547616914821
    "M24" = 1.081471671859142
    "N24" = 1.025653883297145
    "B25" = 1.02
    "C25" = 1.063130745655262
    "D25" = 1.065594896230092
    "E25" = 1.075771949209038
    "F25" = 1.0806812501128
    "I25" = 1.043185712163715
    "J25" = 1.068563862875653
    "K25" = 1.068556641381309
    "L25" = 1.078703938655867
    "M25" = 1.083599103455468
    "N25" = 1.026337398493371
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Updated" $updates.Count "cells"
